$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C1: header cell "enrollmentNumber" gets a new style (bold, centered,
#     no fill/border) distinct from the other header cells (style 5).
#     Build it by pasting formats from a blank, never-touched cell (so no
#     stray fill/border survives) and then making it bold + centered.
$blank1 = $ws.Range("Z100")
$blank1.Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108

# --- A2: type-hint cell for candidateDetailsId changes from "Integer" to
#     "Any" and picks up the light-blue/bordered style used by the other
#     "Any" cells (copy format from an existing "Any"-styled cell, e.g. H2).
#     The sheet is protected, so briefly unlock the cell to allow the write,
#     then restore its locked state.
$srcA2 = $ws.Range("H2")
$srcA2.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Locked = $false
$ws.Range("A2").Value = "Any"
$ws.Range("A2").Locked = $true

# --- Selection moves from C6 to F8
$ws.Range("F8").Select()

$excel.CutCopyMode = $false
